$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 label and values
$ws.Range("A2").Value = "Error Rate"
$ws.Range("B2").Value = 0.05691056910569103
$ws.Range("C2").Value = 0.02439024390243905

# Remove column D contents (Exp 3 header and its value) so the used range shrinks to A1:C2
$ws.Columns.Item(4).Delete()
